$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").ClearContents()
Write-Host "A2 after clear: '$($ws2.Range('A2').Value())'"
$ws2.Range("D2").ClearContents()
Write-Host "D2 after clear: '$($ws2.Range('D2').Value())'"
